$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for week 46 and 47 (semana 48 de 2024 update)
$ws.Range("B47").Value = 478
$ws.Range("B48").Value = 474

# Add new rows for week 48 and 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 378

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 3
